$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 30002110
$ws.Range("J32").Value = 11113455
$ws.Range("L32").Value = 11113455
$ws.Range("N32").Value = -11114107
$ws.Range("H41").Value = 420.3
$ws.Range("I41").Value = 325.2
$ws.Range("J41").Value = 515.4
$ws.Range("K41").Value = 325.2
$ws.Range("L41").Value = 515.4
$ws.Range("M41").Value = 114.8
$ws.Range("N41").Value = -1395.4
$ws.Range("H53").Value = 245.68
$ws.Range("I53").Value = 297.53845
$ws.Range("J53").Value = 189.5
$ws.Range("K53").Value = 297.53845
$ws.Range("L53").Value = 189.5
$ws.Range("M53").Value = 339.46155
$ws.Range("N53").Value = -1463.5
$ws.Range("H57").Value = 127954.125
$ws.Range("J57").Value = 127954.125
$ws.Range("L57").Value = 383862.375
$ws.Range("N57").Value = -384860.375
$ws.Range("H70").Value = 1596.5
$ws.Range("I70").Value = 1557.75
$ws.Range("J70").Value = 1712.75
$ws.Range("K70").Value = 4673.25
$ws.Range("L70").Value = 5138.25
$ws.Range("M70").Value = -4403.25
$ws.Range("N70").Value = -5678.25
$ws.Range("H73").Value = 1596.5
$ws.Range("I73").Value = 1557.75
$ws.Range("J73").Value = 1712.75
$ws.Range("K73").Value = 4673.25
$ws.Range("L73").Value = 5138.25
$ws.Range("M73").Value = -3737.25
$ws.Range("N73").Value = -7010.25
$ws.Range("H107").Value = 5404.8335
$ws.Range("I107").Value = 3873.3333
$ws.Range("K107").Value = 3873.3333
$ws.Range("M107").Value = -1953.3333
$ws.Range("H137").Value = 3000.0334
$ws.Range("J137").Value = 3104
$ws.Range("L137").Value = 9312
$ws.Range("N137").Value = -14412

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 952.5
$ws.Range("I3").Value = 952.5
$ws.Range("K3").Value = 952.5
$ws.Range("M3").Value = -837.5
$ws.Range("H4").Value = 452.66666
$ws.Range("J4").Value = 1122.5
$ws.Range("L4").Value = 1122.5
$ws.Range("N4").Value = -1354.5
$ws.Range("H61").Value = 2654.4614
$ws.Range("I61").Value = 1522.6562
$ws.Range("J61").Value = 7828.4287
$ws.Range("K61").Value = 1522.6562
$ws.Range("L61").Value = 7828.4287
$ws.Range("M61").Value = -1310.6562
$ws.Range("N61").Value = -8252.4287
$ws.Range("H74").Value = 5749.5
$ws.Range("I74").Value = 1499
$ws.Range("J74").Value = 7166.3335
$ws.Range("K74").Value = 1499
$ws.Range("L74").Value = 7166.3335
$ws.Range("M74").Value = -625
$ws.Range("N74").Value = -8914.333500000001
$ws.Range("H77").Value = 5749.5
$ws.Range("I77").Value = 1499
$ws.Range("J77").Value = 7166.3335
$ws.Range("K77").Value = 7495
$ws.Range("L77").Value = 35831.6675
$ws.Range("M77").Value = -3127
$ws.Range("N77").Value = -44567.6675
$ws.Range("H110").Value = 2341.074
$ws.Range("I110").Value = 2210.0952
$ws.Range("K110").Value = 2210.0952
$ws.Range("M110").Value = -165.0952000000002
$ws.Range("H136").Value = 2654.4614
$ws.Range("I136").Value = 1522.6562
$ws.Range("J136").Value = 7828.4287
$ws.Range("K136").Value = 4567.9686
$ws.Range("L136").Value = 23485.2861
$ws.Range("M136").Value = -2017.9686
$ws.Range("N136").Value = -28585.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1515735.9
$ws.Range("J80").Value = 2222787.2
$ws.Range("L80").Value = 2222787.2
$ws.Range("N80").Value = -2224783.2
$ws.Range("H83").Value = 1515735.9
$ws.Range("J83").Value = 2222787.2
$ws.Range("L83").Value = 11113936
$ws.Range("N83").Value = -11123920
$ws.Range("H134").Value = 3220.5881
$ws.Range("I134").Value = 2472.4614
$ws.Range("J134").Value = 5652
$ws.Range("K134").Value = 7417.3842
$ws.Range("L134").Value = 16956
$ws.Range("M134").Value = -4882.3842
$ws.Range("N134").Value = -22026

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2438.8462
$ws.Range("I31").Value = 1973.3636
$ws.Range("J31").Value = 4999
$ws.Range("K31").Value = 1973.3636
$ws.Range("L31").Value = 4999
$ws.Range("M31").Value = -1678.3636
$ws.Range("N31").Value = -5589
$ws.Range("H34").Value = 2438.8462
$ws.Range("I34").Value = 1973.3636
$ws.Range("J34").Value = 4999
$ws.Range("K34").Value = 1973.3636
$ws.Range("L34").Value = 4999
$ws.Range("M34").Value = -1771.3636
$ws.Range("N34").Value = -5403
$ws.Range("H50").Value = 46709
$ws.Range("H51").Value = 35453.555
$ws.Range("J51").Value = 57498
$ws.Range("L51").Value = 57498
$ws.Range("N51").Value = -58970
$ws.Range("H59").Value = 41315.715
$ws.Range("I59").Value = 17523.75
$ws.Range("J59").Value = 73038.336
$ws.Range("K59").Value = 17523.75
$ws.Range("L59").Value = 73038.336
$ws.Range("M59").Value = -16378.75
$ws.Range("N59").Value = -75328.336
$ws.Range("H60").Value = 54462.727
$ws.Range("I60").Value = 22727.285
$ws.Range("J60").Value = 109999.75
$ws.Range("K60").Value = 22727.285
$ws.Range("L60").Value = 109999.75
$ws.Range("M60").Value = -22216.285
$ws.Range("N60").Value = -111021.75
$ws.Range("H61").Value = 35453.555
$ws.Range("J61").Value = 57498
$ws.Range("L61").Value = 57498
$ws.Range("N61").Value = -58194
$ws.Range("H94").Value = 14996.4
$ws.Range("I94").Value = 19667.334
$ws.Range("K94").Value = 19667.334
$ws.Range("M94").Value = -19216.334
$ws.Range("H132").Value = 6882
$ws.Range("I132").Value = 8639.049999999999
$ws.Range("J132").Value = 2977.4443
$ws.Range("K132").Value = 25917.15
$ws.Range("L132").Value = 8932.332900000001
$ws.Range("M132").Value = -23387.15
$ws.Range("N132").Value = -13992.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value = 50293.668
$ws.Range("J115").Value = 50293.668
$ws.Range("L115").Value = 50293.668
$ws.Range("N115").Value = -52643.668
$ws.Range("H136").Value = 6214541
$ws.Range("I136").Value = 12003655
$ws.Range("J136").Value = 11918.429
$ws.Range("K136").Value = 36010965
$ws.Range("L136").Value = 35755.287
$ws.Range("M136").Value = -36008415
$ws.Range("N136").Value = -40855.287
$ws.Range("H139").Value = 94199
$ws.Range("J139").Value = 94199
$ws.Range("L139").Value = 94199
$ws.Range("N139").Value = -104479

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 16148.5
$ws.Range("I38").Value = 12399
$ws.Range("J38").Value = 19898
$ws.Range("K38").Value = 12399
$ws.Range("L38").Value = 19898
$ws.Range("M38").Value = -11926
$ws.Range("N38").Value = -20844
$ws.Range("H62").Value = 6939.5
$ws.Range("I62").Value = 8065.3335
$ws.Range("K62").Value = 8065.3335
$ws.Range("M62").Value = -7441.3335
$ws.Range("H65").Value = 6939.5
$ws.Range("I65").Value = 8065.3335
$ws.Range("K65").Value = 40326.6675
$ws.Range("M65").Value = -37206.6675
$ws.Range("H122").Value = 12843.417
$ws.Range("I122").Value = 7697.4736
$ws.Range("K122").Value = 23092.4208
$ws.Range("M122").Value = -20642.4208
$ws.Range("H136").Value = 3270.6365
$ws.Range("I136").Value = 1974.0769
$ws.Range("J136").Value = 5143.4443
$ws.Range("K136").Value = 5922.2307
$ws.Range("L136").Value = 15430.3329
$ws.Range("M136").Value = -3372.2307
$ws.Range("N136").Value = -20530.3329
$ws.Range("H139").Value = 69804.86
$ws.Range("J139").Value = 69939
$ws.Range("L139").Value = 69939
$ws.Range("N139").Value = -80219
